$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (ALC)
$ws.Range("H6").Value = 279.29413
$ws.Range("I6").Value = 189
$ws.Range("K6").Value = 567
$ws.Range("M6").Value = -455

# Row 53 (ALC)
$ws.Range("H53").Value = 25050306
$ws.Range("I53").Value = 62625108
$ws.Range("J53").Value = 438.5
$ws.Range("K53").Value = 62625108
$ws.Range("L53").Value = 438.5
$ws.Range("M53").Value = -62624471
$ws.Range("N53").Value = -1712.5

# Row 54 (ALC)
$ws.Range("H54").Value = 12000
$ws.Range("J54").Value = 12000
$ws.Range("L54").Value = 12000
$ws.Range("N54").Value = -12972

# Row 55 (ALC)
$ws.Range("H55").Value = 686.4
$ws.Range("I55").Value = 986.6667
$ws.Range("K55").Value = 986.6667
$ws.Range("M55").Value = -772.6667

# Row 129 (ALC)
$ws.Range("H129").Value = 1039.1321
$ws.Range("I129").Value = 547.625
$ws.Range("J129").Value = 1126.5111
$ws.Range("K129").Value = 1642.875
$ws.Range("L129").Value = 3379.5333
$ws.Range("M129").Value = 3357.125
$ws.Range("N129").Value = -13379.5333

# Row 132 (ALC)
$ws.Range("H132").Value = 1245.0465
$ws.Range("I132").Value = 1263.425
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3790.275
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -1260.275
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("ARM")
# Row 53 (ARM)
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

# Row 54 (ARM)
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").ClearContents()
$ws.Range("N54").Value = 0

# Row 61 (ARM)
$ws.Range("H61").Value = 2943.1794
$ws.Range("I61").Value = 2036.5
$ws.Range("J61").Value = 3897.5789
$ws.Range("K61").Value = 2036.5
$ws.Range("L61").Value = 3897.5789
$ws.Range("M61").Value = -1824.5
$ws.Range("N61").Value = -4321.5789

# Row 136 (ARM)
$ws.Range("H136").Value = 2943.1794
$ws.Range("I136").Value = 2036.5
$ws.Range("J136").Value = 3897.5789
$ws.Range("K136").Value = 6109.5
$ws.Range("L136").Value = 11692.7367
$ws.Range("M136").Value = -3559.5
$ws.Range("N136").Value = -16792.7367

$ws = $wb.Worksheets.Item("CRP")
# Row 15 (CRP)
$ws.Range("H15").Value = 2377.25
$ws.Range("J15").Value = 2669.6667
$ws.Range("L15").Value = 2669.6667
$ws.Range("N15").Value = -3009.6667

# Row 31 (CRP)
$ws.Range("H31").Value = 6948.077
$ws.Range("J31").Value = 6948.077
$ws.Range("L31").Value = 6948.077
$ws.Range("N31").Value = -7538.077

# Row 34 (CRP)
$ws.Range("H34").Value = 6948.077
$ws.Range("J34").Value = 6948.077
$ws.Range("L34").Value = 6948.077
$ws.Range("N34").Value = -7352.077

# Row 51 (CRP)
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").ClearContents()
$ws.Range("N51").Value = 0

# Row 58 (CRP)
$ws.Range("H58").Value = 1380.8334
$ws.Range("I58").Value = 1183.7941
$ws.Range("J58").Value = 1715.8
$ws.Range("K58").Value = 1183.7941
$ws.Range("L58").Value = 1715.8
$ws.Range("M58").Value = -980.7941000000001
$ws.Range("N58").Value = -2121.8

# Row 59 (CRP)
$ws.Range("H59").Value = 14039.7
$ws.Range("I59").Value = 3000
$ws.Range("J59").Value = 15266.333
$ws.Range("K59").Value = 3000
$ws.Range("L59").Value = 15266.333
$ws.Range("M59").Value = -1855
$ws.Range("N59").Value = -17556.333

# Row 60 (CRP)
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = 0

# Row 61 (CRP)
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").ClearContents()
$ws.Range("N61").Value = 0

# Row 94 (CRP)
$ws.Range("H94").Value = 4600.1035
$ws.Range("I94").Value = 4822
$ws.Range("J94").Value = 4464.5
$ws.Range("K94").Value = 4822
$ws.Range("L94").Value = 4464.5
$ws.Range("M94").Value = -4371
$ws.Range("N94").Value = -5366.5

# Row 134 (CRP)
$ws.Range("H134").Value = 3556.2666
$ws.Range("I134").Value = 3611.5
$ws.Range("J134").Value = 3335.3333
$ws.Range("K134").Value = 10834.5
$ws.Range("L134").Value = 10005.9999
$ws.Range("M134").Value = -8299.5
$ws.Range("N134").Value = -15075.9999

# Row 136 (CRP)
$ws.Range("H136").Value = 1380.8334
$ws.Range("I136").Value = 1183.7941
$ws.Range("J136").Value = 1715.8
$ws.Range("K136").Value = 3551.3823
$ws.Range("L136").Value = 5147.4
$ws.Range("M136").Value = -1001.3823
$ws.Range("N136").Value = -10247.4

$ws = $wb.Worksheets.Item("CUL")
# Row 25 (CUL)
$ws.Range("H25").Value = 2750
$ws.Range("I25").Value = 500
$ws.Range("J25").Value = 5000
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = -1331
$ws.Range("N25").Value = -15338

# Row 30 (CUL)
$ws.Range("H30").Value = 2750
$ws.Range("I30").Value = 500
$ws.Range("J30").Value = 5000
$ws.Range("K30").Value = 1500
$ws.Range("L30").Value = 15000
$ws.Range("M30").Value = -1398
$ws.Range("N30").Value = -15204

# Row 54 (CUL)
$ws.Range("H54").Value = 3132.0454
$ws.Range("J54").Value = 3132.0454
$ws.Range("L54").Value = 9396.136200000001
$ws.Range("N54").Value = -10514.1362

# Row 55 (CUL)
$ws.Range("H55").Value = 3004.5454
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3004.5454
$ws.Range("K55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("M55").Value = 9013.636200000001
$ws.Range("N55").Value = -9367.636200000001

# Row 87 (CUL)
$ws.Range("H87").Value = 5954.6665
$ws.Range("I87").Value = 5954.6665
$ws.Range("K87").Value = 17863.9995
$ws.Range("M87").Value = -16615.9995

# Row 90 (CUL)
$ws.Range("H90").Value = 5954.6665
$ws.Range("I90").Value = 5954.6665
$ws.Range("K90").Value = 53591.9985
$ws.Range("M90").Value = -47351.9985

# Row 92 (CUL)
$ws.Range("H92").Value = 776.6667
$ws.Range("I92").Value = 501.5
$ws.Range("J92").Value = 796.3214
$ws.Range("K92").Value = 1504.5
$ws.Range("L92").Value = 2388.9642
$ws.Range("M92").Value = -256.5
$ws.Range("N92").Value = -4884.9642

# Row 98 (CUL)
$ws.Range("H98").Value = 11111806
$ws.Range("I98").Value = 312.5
$ws.Range("K98").Value = 937.5
$ws.Range("M98").Value = 560.5

# Row 109 (CUL)
$ws.Range("H109").Value = 2448
$ws.Range("I109").Value = 2000
$ws.Range("J109").Value = 2497.7778
$ws.Range("K109").Value = 6000
$ws.Range("L109").Value = 7493.3334
$ws.Range("M109").Value = -4960
$ws.Range("N109").Value = -9573.3334

# Row 137 (CUL)
$ws.Range("H137").Value = 9269.972
$ws.Range("I137").Value = 7007.5
$ws.Range("J137").Value = 12286.6
$ws.Range("K137").Value = 21022.5
$ws.Range("L137").Value = 36859.8
$ws.Range("M137").Value = -15922.5
$ws.Range("N137").Value = -47059.8

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM)
$ws.Range("H70").Value = 5633.5454
$ws.Range("J70").Value = 5409.4287
$ws.Range("L70").Value = 5409.4287
$ws.Range("N70").Value = -5949.4287

# Row 73 (GSM)
$ws.Range("H73").Value = 5633.5454
$ws.Range("J73").Value = 5409.4287
$ws.Range("L73").Value = 5409.4287
$ws.Range("N73").Value = -7281.4287

# Row 80 (GSM)
$ws.Range("H80").Value = 2469.4119
$ws.Range("I80").Value = 2456.6667
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 2456.6667
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -1458.6667
$ws.Range("N80").Value = -4496

# Row 83 (GSM)
$ws.Range("H83").Value = 2469.4119
$ws.Range("I83").Value = 2456.6667
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 12283.3335
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -7291.333500000001
$ws.Range("N83").Value = -22484

# Row 126 (GSM)
$ws.Range("H126").Value = 6124.48
$ws.Range("I126").Value = 7367.3335
$ws.Range("J126").Value = 2928.5715
$ws.Range("K126").Value = 22102.0005
$ws.Range("L126").Value = 8785.7145
$ws.Range("M126").Value = -19632.0005
$ws.Range("N126").Value = -13725.7145

# Row 132 (GSM)
$ws.Range("H132").Value = 5441.439
$ws.Range("I132").Value = 10594.733
$ws.Range("J132").Value = 2468.3845
$ws.Range("K132").Value = 31784.199
$ws.Range("L132").Value = 7405.1535
$ws.Range("M132").Value = -29254.199
$ws.Range("N132").Value = -12465.1535

$ws = $wb.Worksheets.Item("LTW")
# Row 55 (LTW)
$ws.Range("H55").Value = 350.2353
$ws.Range("I55").Value = 318.63635
$ws.Range("J55").Value = 408.16666
$ws.Range("K55").Value = 318.63635
$ws.Range("L55").Value = 408.16666
$ws.Range("M55").Value = -145.63635
$ws.Range("N55").Value = -754.16666

# Row 82 (LTW)
$ws.Range("H82").Value = 17016598
$ws.Range("I82").Value = 2500947.5
$ws.Range("K82").Value = 2500947.5
$ws.Range("M82").Value = -2500586.5

# Row 85 (LTW)
$ws.Range("H85").Value = 17016598
$ws.Range("I85").Value = 2500947.5
$ws.Range("K85").Value = 2500947.5
$ws.Range("M85").Value = -2499699.5

# Row 132 (LTW)
$ws.Range("H132").Value = 10103226
$ws.Range("I132").Value = 14494724
$ws.Range("J132").Value = 2780.4
$ws.Range("K132").Value = 43484172
$ws.Range("L132").Value = 8341.200000000001
$ws.Range("M132").Value = -43481642
$ws.Range("N132").Value = -13401.2

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (WVR)
$ws.Range("H122").Value = 1846.875
$ws.Range("I122").Value = 1595.8334
$ws.Range("K122").Value = 4787.5002
$ws.Range("M122").Value = -2337.5002
